$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 34, pushing current rows 34 and 35 down to 35 and 36
$ws.Rows.Item(34).Insert()

# Populate the new row 34 with the updated record (copy of old row 34 data, with
# date / volume / price updates applied per the diff)
$ws.Cells.Item(34, 1).Value = 3
$ws.Cells.Item(34, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(34, 3).Value = "Coquimbo"
$ws.Cells.Item(34, 4).Value = 44714
$ws.Cells.Item(34, 4).NumberFormat = $ws.Cells.Item(35, 4).NumberFormat
$ws.Cells.Item(34, 5).Value = 5
$ws.Cells.Item(34, 6).Value = 100112035
$ws.Cells.Item(34, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(34, 8).Value = "Sin especificar"
$ws.Cells.Item(34, 9).Value = "Primera"
$ws.Cells.Item(34, 10).Value = 100
$ws.Cells.Item(34, 11).Value = 15000
$ws.Cells.Item(34, 12).Value = 15500
$ws.Cells.Item(34, 13).Value = 15250
$ws.Cells.Item(34, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(34, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(34, 16).Value = 1017
$ws.Cells.Item(34, 17).Value = 15
$ws.Cells.Item(34, 18).Value = "Hortaliza"
